$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Row 26: new backlog activity
$ws.Range("A26").Value = (Get-Date -Year 2024 -Month 5 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B26").Value = "Patryk"
$ws.Range("C26").Value = "take_photo: implementacja przycisku usuwania"
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = "Ukończono"

# Row 27: new backlog activity
$ws.Range("A27").Value = (Get-Date -Year 2024 -Month 5 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B27").Value = "Patryk"
$ws.Range("C27").Value = "take_photo: implementacja przycisków wyjścia i zapisu"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 35
$ws.Range("F27").Value = "Ukończono"

# Scroll/selection change seen in the diff (new active cell + top row of view)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("I31").Select()
